# Add new rows of "bombardier" load-test targets (rows 72-86 on Sheet1).
# For each row we set column B (IP address) and column C (company/tag),
# then put the "sudo docker run ..." formula in column A that references
# the same row's B cell - matching the existing pattern used by every
# other populated row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @{ Row = 72; Ip = "213.129.114.94";  Tag = "efko" },
    @{ Row = 73; Ip = "213.129.114.95";  Tag = "efko" },
    @{ Row = 74; Ip = "195.239.33.178";  Tag = "efko" },
    @{ Row = 75; Ip = "213.129.114.91";  Tag = "efko" },
    @{ Row = 76; Ip = "213.129.114.92";  Tag = "efko" },
    @{ Row = 77; Ip = "213.129.115.54";  Tag = "efko" },
    @{ Row = 78; Ip = "213.129.114.93";  Tag = "efko" },
    @{ Row = 79; Ip = "213.129.114.89";  Tag = "efko" },
    @{ Row = 80; Ip = "213.129.114.95";  Tag = "efko" },
    @{ Row = 81; Ip = "213.129.115.50";  Tag = "efko" },
    @{ Row = 82; Ip = "213.129.114.88";  Tag = "efko" },
    @{ Row = 83; Ip = "195.239.33.180";  Tag = "efko" },
    @{ Row = 84; Ip = "195.43.90.110";   Tag = "homecredit" },
    @{ Row = 85; Ip = "178.154.197.231"; Tag = "samolet_ru" },
    @{ Row = 86; Ip = "82.202.246.121";  Tag = "samolet_ru" }
)

foreach ($entry in $data) {
    $r = $entry.Row

    $ws.Cells.Item($r, 2).Value = $entry.Ip
    $ws.Cells.Item($r, 3).Value = $entry.Tag
    $ws.Cells.Item($r, 1).Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B' + $r + '&"&& sleep 5;"'
}
